# Fruta / hortaliza, semanal
# The weekly data window rolled forward by one record: the newest record
# (old row 5) becomes the first data row (row 2), and the remaining rows
# shift down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the "before" values for the columns that move (D, L, M, N, O, P, S)
# across data rows 2-5 before overwriting anything.
$cols = @("D", "L", "M", "N", "O", "P", "S")
$before = @{}
foreach ($row in 2..5) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Range("$col$row").Value2
    }
    $before[$row] = $rowVals
}

# New row 2 = old row 5; new row r = old row (r-1) for r = 3,4,5
$ws.Range("D2").Value = $before[5]["D"]
$ws.Range("L2").Value = $before[5]["L"]
$ws.Range("M2").Value = $before[5]["M"]
$ws.Range("N2").Value = $before[5]["N"]
$ws.Range("O2").Value = $before[5]["O"]
$ws.Range("P2").Value = $before[5]["P"]
$ws.Range("S2").Value = $before[5]["S"]

$ws.Range("D3").Value = $before[2]["D"]
$ws.Range("L3").Value = $before[2]["L"]
$ws.Range("M3").Value = $before[2]["M"]
$ws.Range("N3").Value = $before[2]["N"]
$ws.Range("O3").Value = $before[2]["O"]
$ws.Range("P3").Value = $before[2]["P"]
$ws.Range("S3").Value = $before[2]["S"]

$ws.Range("D4").Value = $before[3]["D"]
$ws.Range("L4").Value = $before[3]["L"]
$ws.Range("M4").Value = $before[3]["M"]
$ws.Range("N4").Value = $before[3]["N"]
$ws.Range("O4").Value = $before[3]["O"]
$ws.Range("P4").Value = $before[3]["P"]
$ws.Range("S4").Value = $before[3]["S"]

$ws.Range("D5").Value = $before[4]["D"]
$ws.Range("L5").Value = $before[4]["L"]
$ws.Range("M5").Value = $before[4]["M"]
$ws.Range("N5").Value = $before[4]["N"]
$ws.Range("O5").Value = $before[4]["O"]
$ws.Range("P5").Value = $before[4]["P"]
$ws.Range("S5").Value = $before[4]["S"]
